$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 147
$ws.Cells.Item(147, 2).Value = 5456594
$ws.Cells.Item(147, 6).Value = "Rakow Czestochowa"
$ws.Cells.Item(147, 7).Value = "Zaglebie Lubin"
$ws.Cells.Item(147, 8).Value = 1
$ws.Cells.Item(147, 9).Value = 1
$ws.Cells.Item(147, 11).Value = 1.444
$ws.Cells.Item(147, 12).Value = 4.5
$ws.Cells.Item(147, 13).Value = 5.75
$ws.Cells.Item(147, 14).Value = 1.3
$ws.Cells.Item(147, 15).Value = 5.25
$ws.Cells.Item(147, 16).Value = 7
$ws.Cells.Item(147, 17).Value = -1.5
$ws.Cells.Item(147, 18).Value = 1.9
$ws.Cells.Item(147, 19).Value = 1.95
$ws.Cells.Item(147, 20).Value = 3
$ws.Cells.Item(147, 21).Value = 1.9
$ws.Cells.Item(147, 22).Value = 1.95
$ws.Cells.Item(147, 24).Value = 4.25
$ws.Cells.Item(147, 26).Value = -1
$ws.Cells.Item(147, 27).Value = 0.95
$ws.Cells.Item(147, 29).Value = 0.95
# Row 148
$ws.Cells.Item(148, 2).Value = 5456603
$ws.Cells.Item(148, 6).Value = "Lech Poznan"
$ws.Cells.Item(148, 7).Value = "Jagiellonia Bialystok"
$ws.Cells.Item(148, 8).Value = 2
$ws.Cells.Item(148, 11).Value = 1.363
$ws.Cells.Item(148, 12).Value = 4.75
$ws.Cells.Item(148, 13).Value = 6.5
$ws.Cells.Item(148, 14).Value = 1.222
$ws.Cells.Item(148, 15).Value = 5.5
$ws.Cells.Item(148, 16).Value = 8
$ws.Cells.Item(148, 17).Value = -1.75
$ws.Cells.Item(148, 18).Value = 1.925
$ws.Cells.Item(148, 19).Value = 1.925
$ws.Cells.Item(148, 20).Value = 3.25
$ws.Cells.Item(148, 21).Value = 1.95
$ws.Cells.Item(148, 22).Value = 1.9
$ws.Cells.Item(148, 23).Value = 0.222
$ws.Cells.Item(148, 26).Value = 0.4625
$ws.Cells.Item(148, 27).Value = -0.5
$ws.Cells.Item(148, 28).Value = -1
$ws.Cells.Item(148, 29).Value = 0.8999999999999999
# Row 149
$ws.Cells.Item(149, 2).Value = 5460884
$ws.Cells.Item(149, 6).Value = "Miedz Legnica"
$ws.Cells.Item(149, 7).Value = "Gornik Zabrze"
$ws.Cells.Item(149, 9).Value = 0
$ws.Cells.Item(149, 10).Value = "D"
$ws.Cells.Item(149, 11).Value = 3.6
$ws.Cells.Item(149, 12).Value = 3.5
$ws.Cells.Item(149, 13).Value = 1.909
$ws.Cells.Item(149, 14).Value = 3
$ws.Cells.Item(149, 15).Value = 3.5
$ws.Cells.Item(149, 16).Value = 2.1
$ws.Cells.Item(149, 17).Value = 0.25
$ws.Cells.Item(149, 18).Value = 1.95
$ws.Cells.Item(149, 19).Value = 1.9
$ws.Cells.Item(149, 20).Value = 2.75
$ws.Cells.Item(149, 21).Value = 1.975
$ws.Cells.Item(149, 22).Value = 1.875
$ws.Cells.Item(149, 24).Value = 2.5
$ws.Cells.Item(149, 25).Value = -1
$ws.Cells.Item(149, 26).Value = 0.475
$ws.Cells.Item(149, 27).Value = -0.5
$ws.Cells.Item(149, 28).Value = -1
$ws.Cells.Item(149, 29).Value = 0.875
# Row 151
$ws.Cells.Item(151, 2).Value = 5465446
$ws.Cells.Item(151, 6).Value = "Cracovia Krakow"
$ws.Cells.Item(151, 7).Value = "Wisla Plock"
$ws.Cells.Item(151, 8).Value = 3
$ws.Cells.Item(151, 11).Value = 2.15
$ws.Cells.Item(151, 12).Value = 3.5
$ws.Cells.Item(151, 13).Value = 2.875
$ws.Cells.Item(151, 14).Value = 2.25
$ws.Cells.Item(151, 15).Value = 3.6
$ws.Cells.Item(151, 16).Value = 2.7
$ws.Cells.Item(151, 17).Value = -0.25
$ws.Cells.Item(151, 18).Value = 2.05
$ws.Cells.Item(151, 19).Value = 1.75
$ws.Cells.Item(151, 20).Value = 2.5
$ws.Cells.Item(151, 21).Value = 1.825
$ws.Cells.Item(151, 22).Value = 2.025
$ws.Cells.Item(151, 23).Value = 1.25
$ws.Cells.Item(151, 26).Value = 1.05
$ws.Cells.Item(151, 28).Value = 0.825
$ws.Cells.Item(151, 29).Value = -1
# Row 152
$ws.Cells.Item(152, 2).Value = 5428774
$ws.Cells.Item(152, 6).Value = "Pogon Szczecin"
$ws.Cells.Item(152, 7).Value = "Radomiak Radom"
$ws.Cells.Item(152, 8).Value = 4
$ws.Cells.Item(152, 11).Value = 1.571
$ws.Cells.Item(152, 12).Value = 4
$ws.Cells.Item(152, 13).Value = 4.75
$ws.Cells.Item(152, 14).Value = 1.533
$ws.Cells.Item(152, 15).Value = 4.333
$ws.Cells.Item(152, 16).Value = 4.75
$ws.Cells.Item(152, 17).Value = -1
$ws.Cells.Item(152, 18).Value = 1.875
$ws.Cells.Item(152, 19).Value = 1.975
$ws.Cells.Item(152, 20).Value = 3
$ws.Cells.Item(152, 21).Value = 1.875
$ws.Cells.Item(152, 22).Value = 1.975
$ws.Cells.Item(152, 23).Value = 0.5329999999999999
$ws.Cells.Item(152, 26).Value = 0.875
$ws.Cells.Item(152, 27).Value = -1
$ws.Cells.Item(152, 28).Value = 0.875
$ws.Cells.Item(152, 29).Value = -1
# Row 153
$ws.Cells.Item(153, 2).Value = 5467427
$ws.Cells.Item(153, 6).Value = "Stal Mielec"
$ws.Cells.Item(153, 7).Value = "Warta Poznan"
$ws.Cells.Item(153, 9).Value = 0
$ws.Cells.Item(153, 10).Value = "H"
$ws.Cells.Item(153, 11).Value = 2.375
$ws.Cells.Item(153, 12).Value = 3.2
$ws.Cells.Item(153, 13).Value = 2.8
$ws.Cells.Item(153, 14).Value = 2.6
$ws.Cells.Item(153, 15).Value = 3.1
$ws.Cells.Item(153, 16).Value = 2.625
$ws.Cells.Item(153, 17).Value = 0
$ws.Cells.Item(153, 18).Value = 1.925
$ws.Cells.Item(153, 19).Value = 1.925
$ws.Cells.Item(153, 20).Value = 2.25
$ws.Cells.Item(153, 21).Value = 1.975
$ws.Cells.Item(153, 22).Value = 1.875
$ws.Cells.Item(153, 23).Value = 1.6
$ws.Cells.Item(153, 24).Value = -1
$ws.Cells.Item(153, 26).Value = 0.925
$ws.Cells.Item(153, 27).Value = -1
$ws.Cells.Item(153, 29).Value = 0.875
# Row 154
$ws.Cells.Item(154, 2).Value = 5461475
$ws.Cells.Item(154, 6).Value = "Widzew Lodz"
$ws.Cells.Item(154, 7).Value = "Korona Kielce"
$ws.Cells.Item(154, 8).Value = 0
$ws.Cells.Item(154, 9).Value = 3
$ws.Cells.Item(154, 10).Value = "A"
$ws.Cells.Item(154, 11).Value = 2.1
$ws.Cells.Item(154, 12).Value = 3.3
$ws.Cells.Item(154, 13).Value = 3.2
$ws.Cells.Item(154, 14).Value = 2.375
$ws.Cells.Item(154, 15).Value = 3.3
$ws.Cells.Item(154, 16).Value = 2.7
$ws.Cells.Item(154, 17).Value = 0
$ws.Cells.Item(154, 18).Value = 1.8
$ws.Cells.Item(154, 19).Value = 2.05
$ws.Cells.Item(154, 20).Value = 2.5
$ws.Cells.Item(154, 21).Value = 1.825
$ws.Cells.Item(154, 22).Value = 2.025
$ws.Cells.Item(154, 23).Value = -1
$ws.Cells.Item(154, 25).Value = 1.7
$ws.Cells.Item(154, 26).Value = -1
$ws.Cells.Item(154, 27).Value = 1.05
$ws.Cells.Item(154, 28).Value = 0.825
# Row 341
$ws.Cells.Item(341, 8).Value = 0
$ws.Cells.Item(341, 9).Value = 0
$ws.Cells.Item(341, 10).Value = "D"
$ws.Cells.Item(341, 14).Value = 2.3
$ws.Cells.Item(341, 15).Value = 2.875
$ws.Cells.Item(341, 16).Value = 3.5
$ws.Cells.Item(341, 18).Value = 1.925
$ws.Cells.Item(341, 19).Value = 1.925
$ws.Cells.Item(341, 21).Value = 1.925
$ws.Cells.Item(341, 22).Value = 1.925
$ws.Cells.Item(341, 23).Value = -1
$ws.Cells.Item(341, 24).Value = 1.875
$ws.Cells.Item(341, 25).Value = -1
$ws.Cells.Item(341, 26).Value = -0.5
$ws.Cells.Item(341, 27).Value = 0.4625
$ws.Cells.Item(341, 28).Value = -1
$ws.Cells.Item(341, 29).Value = 0.925
# Row 342
$ws.Cells.Item(342, 8).Value = 4
$ws.Cells.Item(342, 9).Value = 2
$ws.Cells.Item(342, 10).Value = "H"
$ws.Cells.Item(342, 14).Value = 1.4
$ws.Cells.Item(342, 15).Value = 4.75
$ws.Cells.Item(342, 16).Value = 7
$ws.Cells.Item(342, 18).Value = 1.975
$ws.Cells.Item(342, 19).Value = 1.875
$ws.Cells.Item(342, 20).Value = 3.25
$ws.Cells.Item(342, 21).Value = 2.05
$ws.Cells.Item(342, 22).Value = 1.8
$ws.Cells.Item(342, 23).Value = 0.3999999999999999
$ws.Cells.Item(342, 24).Value = -1
$ws.Cells.Item(342, 25).Value = -1
$ws.Cells.Item(342, 26).Value = 0.9750000000000001
$ws.Cells.Item(342, 27).Value = -1
$ws.Cells.Item(342, 28).Value = 1.05
$ws.Cells.Item(342, 29).Value = -1
# Row 343
$ws.Cells.Item(343, 18).Value = 1.975
$ws.Cells.Item(343, 19).Value = 1.875
# Row 344
$ws.Cells.Item(344, 21).Value = 1.875
$ws.Cells.Item(344, 22).Value = 1.975
# Row 345
$ws.Cells.Item(345, 14).Value = 5.25
$ws.Cells.Item(345, 16).Value = 1.615
$ws.Cells.Item(345, 18).Value = 2.05
$ws.Cells.Item(345, 19).Value = 1.8
# Row 346
$ws.Cells.Item(346, 21).Value = 1.95
$ws.Cells.Item(346, 22).Value = 1.9
# Row 347
$ws.Cells.Item(347, 18).Value = 2
$ws.Cells.Item(347, 19).Value = 1.85
$ws.Cells.Item(347, 21).Value = 1.975
$ws.Cells.Item(347, 22).Value = 1.875
# Row 349
$ws.Cells.Item(349, 14).Value = 2.5
$ws.Cells.Item(349, 15).Value = 3.2
$ws.Cells.Item(349, 16).Value = 2.8
$ws.Cells.Item(349, 21).Value = 1.975
$ws.Cells.Item(349, 22).Value = 1.875
